$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = "1:0"
$ws.Range("H8").Value = "0:-1"
$ws.Range("A10").Value = "Saturday, 24 May, 2025 9:32 AM"
